# Adds a second worksheet ("Sheet2") containing a distinctive collexeme
# analysis overview table (mirrors the layout of Sheet1, but for the
# "hangover" example comparing the "mother of all" cxn and the
# "ADJ-est N ever" construction), widens a column on Sheet1, and updates
# the active sheet / selections to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 tweak: new third column width ---------------------------------
$ws1.Columns.Item(3).ColumnWidth = 38.6

# --- Add Sheet2 right after Sheet1 -----------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Column widths for Sheet2 (col B / col C).
$ws2.Columns.Item(2).ColumnWidth = 47.92
$ws2.Columns.Item(3).ColumnWidth = 44.92

# --- Sheet2 content ---------------------------------------------------------
# Header row.
$ws2.Cells.Item(1,2).Value = "Word l" + [char]0x1D62 + " of Class L"
$ws2.Cells.Item(1,3).Value = "Other Words of Class L"
$ws2.Cells.Item(1,4).Value = "Total"
$ws2.Cells.Item(2,1).Value = "Construction c" + [char]0x2081 + " of Class C"

# Body cells (order chosen to match the source workbook's shared-string table).
$ws2.Cells.Item(2,3).Value = 'Frequency of all other nouns in the "mother of all" cxn'
$ws2.Cells.Item(3,2).Value = 'Frequency of "hangover" in the "ADJ-est N ever" construction'
$ws2.Cells.Item(2,2).Value = 'Frequency of "hangover" in the "mother of all" cxn'
$ws2.Cells.Item(3,3).Value = 'Frequency of all other nouns in the "ADJ-est N ever" construction'
$ws2.Cells.Item(4,2).Value = 'Total frequency of "hangover" in both constructions'
$ws2.Cells.Item(4,3).Value = 'Total frequency of all other nouns in the two constructions'
$ws2.Cells.Item(2,4).Value = 'Total frequency of "mother of all"'
$ws2.Cells.Item(3,4).Value = 'Total frequency of "ADJ-est N ever"'
$ws2.Cells.Item(4,4).Value = 'Total frequency of both cxns'

$ws2.Cells.Item(3,1).Value = "Construction c" + [char]0x2082 + " of Class C"
$ws2.Cells.Item(4,1).Value = "Total"

# --- Formatting: reuse Sheet1's existing "Helvetica" 9pt style (style index
# 1) by copying its format only, so no redundant style entries are created.
$ws1.Range("A2").Copy()
$ws2.Range("A1:B4").PasteSpecial(-4122)
$ws2.Range("C2:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selections / active sheet ----------------------------------------------
$ws1.Range("A1:D4").Select()

$ws2.Activate()
$ws2.Range("D5").Select()
